# Automatische test-sync: 2025-06-22 19:11:50
# Adds the new mail-log entry (row 41) to the "Logs" sheet and bumps the
# "IT / Technisch probleem" tally on the "Dashboard" sheet from 6 to 7.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(41, 1).Value = "Kan mijn wachtwoord niet resetten"
$logs.Cells.Item(41, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(41, 3).Value = "Ik krijg geen e-mail bij wachtwoord resetten."
$logs.Cells.Item(41, 4).Value = "IT / Technisch probleem"
$logs.Cells.Item(41, 5).Value = "Geachte klant,`nBedankt voor uw e-mail. Om het probleem met het niet ontvangen van e-mails voor het resetten van uw wachtwoord op te lossen, hebben we wat meer informatie nodig. Zou u ons alstublieft de gebruikersnaam of het e-mailadres kunnen geven waarvoor u het wachtwoord wilt resetten? Op deze manier kunnen we verder onderzoeken waar het probleem precies ligt en u van dienst zijn.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item(41, 6).Value = "2025-06-22 19:11:15"
$logs.Cells.Item(41, 7).Value = "Ja"

# --- Extend the conditional-formatting ranges to cover the new row ---
$catRange = $logs.Range("D2:D41")
$catRange.FormatConditions.Item(1).ModifyAppliesToRange($catRange)

$answeredRange = $logs.Range("G2:G41")
$answeredRange.FormatConditions.Item(1).ModifyAppliesToRange($answeredRange)

# --- Dashboard sheet: "IT / Technisch probleem" count 6 -> 7 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 7
